$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 17.75
$ws.Range("C3").Value = 12.93879281845967
$ws.Range("D3").Value = 12.7541939361204
$ws.Range("E3").Value = 1

# Row 4 (std)
$ws.Range("B4").Value = 0.4578165130223675
$ws.Range("C4").Value = 0.4999456601970534
$ws.Range("D4").Value = 0.4469435261075588
$ws.Range("E4").Value = 0

# Row 5 (min)
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 12.03676160791397
$ws.Range("D5").Value = 12.00034100984755
$ws.Range("E5").Value = 1

# Row 6 (25%)
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 12.20948416866571
$ws.Range("D6").Value = 12.01070529872066
$ws.Range("E6").Value = 1

# Row 7 (50%)
$ws.Range("B7").Value = 18
$ws.Range("C7").Value = 13.18149858061092
$ws.Range("D7").Value = 13.01004228867749
$ws.Range("E7").Value = 1

# Row 8 (75%)
$ws.Range("B8").Value = 18
$ws.Range("C8").Value = 13.25207383865541
$ws.Range("D8").Value = 13.0208299124856
$ws.Range("E8").Value = 1

# Row 9 (max)
$ws.Range("B9").Value = 19
$ws.Range("C9").Value = 13.5611755592543
$ws.Range("D9").Value = 13.06674066044482
$ws.Range("E9").Value = 1
